$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto symbol list snapshot (GitHub Actions scrape).
# Rows 7-24 shift down by one slot because a new coin (GateToken) is
# inserted at row 7; every other listed row only gets fresh Price /
# Volume(1h) figures (columns D and E).
$rows = @(
    @{Row=2;  D='328.78';  E='0.44%'}
    @{Row=3;  D='43.99';   E='0.14%'}
    @{Row=4;  D='5.585';   E='1.63%'}
    @{Row=5;  D='0.08062'; E='-0.08%'}
    @{Row=6;  D='1.998';   E='6.23%'}
    @{Row=7;  B='GateToken';                          C='https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt';                          D='4.334';        E='1.33%'}
    @{Row=8;  B='MXToken';                            C='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx';                            D='0.9514';       E='1.60%'}
    @{Row=9;  B='BTSEToken';                          C='https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse';                        D='2.560';        E='-6.01%'}
    @{Row=10; B='LiechtensteinCryptoassetsExchange';  C='https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx';    D='0.1167';       E='0.54%'}
    @{Row=11; B='WazirX';                             C='https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx';                                D='0.1851';       E='-1.93%'}
    @{Row=12; B='MCDex';                              C='https://coinranking.com/coin/3nMM61qeg+mcdex-mcb';                                D='11.85';        E='38.39%'}
    @{Row=13; B='MandalaExchangeToken';               C='https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx';             D='0.09737';      E='1.86%'}
    @{Row=14; B='BitrueCoin';                         C='https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr';                           D='0.04705';      E='14.52%'}
    @{Row=15; B='BitMartToken';                       C='https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx';                     D='0.1068';       E='0.18%'}
    @{Row=16; B='BitForexToken';                      C='https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf';                     D='0.001292';     E='1.32%'}
    @{Row=17; B='CoinExToken';                        C='https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet';                      D='0.04235';      E='-2.44%'}
    @{Row=18; B='TigerCash';                          C='https://coinranking.com/coin/6hIn06L2+tigercash-tch';                             D='0.006033';     E='1.09%'}
    @{Row=19; D='3.369';   E='-5.72%'}
    @{Row=20; B='BitpandaEcosystemToken';             C='https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best';               D='0.3475';       E='-0.28%'}
    @{Row=21; B='ProBitToken';                        C='https://coinranking.com/coin/lQP4d6T2+probittoken-prob';                          D='0.1411';       E='3.46%'}
    @{Row=22; B='ZBToken';                            C='https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb';                           D='0.2511';       E='-3.11%'}
    @{Row=23; B='BitKan';                             C='https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan';                           D='0.001253';     E='1.61%'}
    @{Row=24; B='HotbitToken';                        C='https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb';                       D='0.004317';     E='-1.76%'}
    @{Row=25; D='0.0001192';    E='-3.32%'}
    @{Row=26; E='-0.48%'}
    @{Row=38; D='0.02607';      E='-1.24%'}
    @{Row=39; D='0.05517';      E='1.54%'}
    @{Row=40; D='0.007565';     E='-0.96%'}
    @{Row=41; D='0.1402';       E='0.94%'}
    @{Row=42; D='0.008093';     E='-29.34%'}
    @{Row=43; D='0.002019';     E='-4.43%'}
    @{Row=44; D='0.008388';     E='-9.37%'}
    @{Row=45; D='0.00007093';   E='2.58%'}
    @{Row=46; E='-0.06%'}
    @{Row=47; E='1.18%'}
    @{Row=48; D='0.004842';     E='35.80%'}
    @{Row=49; D='0.00002104';   E='-0.06%'}
    @{Row=50; E='-0.06%'}
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Coin name / link are plain, non-numeric-looking text -- safe to set directly.
    if ($r.ContainsKey('B')) {
        $ws.Range("B$rowNum").Value = $r.B
    }
    if ($r.ContainsKey('C')) {
        $ws.Range("C$rowNum").Value = $r.C
    }

    # Price / Volume(1h) columns hold numeric- and percent-looking strings
    # that must stay plain text (matching the original inlineStr cells),
    # so force a text number format while writing, then restore the
    # cell's original style afterwards.
    $de = $ws.Range("D$rowNum`:E$rowNum")
    $origStyle = $de.Style
    $de.NumberFormat = "@"
    if ($r.ContainsKey('D')) {
        $ws.Range("D$rowNum").Value = $r.D
    }
    if ($r.ContainsKey('E')) {
        $ws.Range("E$rowNum").Value = $r.E
    }
    $de.Style = $origStyle
}
